$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.213.84"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "1.806.41"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("D5").Value = "'317.10"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").Value = "'0.5319"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").Value = "'0.3775"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "'0.07479"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "'42.02"
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").Value = "'0.9997"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "'6.211"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "'20.55"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "'7.365"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "1.806.75"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "'89.73"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "'0.00001065"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'0.06505"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").Value = "'17.36"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "28.230.80"
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("D24").Value = "'11.21"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'2.084"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("D26").Value = "'155.78"
$ws.Range("E26").Value = "  -2.93%  "
$ws.Range("D27").Value = "'20.47"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "2.014.35"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").Value = "'2.330"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").Value = "'122.12"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("D31").Value = "'1.115"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("D32").Value = "'0.1086"
$ws.Range("E32").Value = "  +6.81%  "
$ws.Range("D33").Value = "'5.583"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").Value = "'3.621"
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").Value = "'0.07137"
$ws.Range("E35").Value = "  +8.75%  "
$ws.Range("D36").Value = "'0.2225"
$ws.Range("E36").Value = "  -2.92%  "
$ws.Range("D37").Value = "'0.02297"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").Value = "'5.084"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").Value = "'8.487"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6179"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'11.14"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.439"
$ws.Range("E42").Value = "  +3.80%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.182"
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("D44").Value = "'13.46"
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").Value = "'3.689"
$ws.Range("D46").Value = "'0.5774"
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("D47").Value = "'125.31"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("D49").Value = "'1.923"
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("D50").Value = "'0.06821"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").Value = "'71.93"
$ws.Range("E51").Value = "  -1.25%  "
